$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-21: update date serial from 45221 to 45224
for ($row = 2; $row -le 21; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45221) {
        $cell.Value2 = 45224
    }
}
